# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the author's fix: the old scraper grabbed team stats but not the
# season record, so here we populate columns AD:AF with the team's
# win/loss/tie totals for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels, same header style as the existing
#     header cells (bold, centered, bordered) by copying format from AC1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-46): same record (86 wins, 76 losses, 0 ties) for the
#     whole roster, since it is one team's season record.
$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 86   # AD
    $ws.Cells.Item($r, 31).Value = 76   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
